# Weekly update: insert 2 new "Sandia" price records at the top of the
# historical block (rows 330-331), shifting all subsequent rows down by 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 330:331 - this shifts existing rows 330..429
# down to 332..431 and copies formatting (e.g. the date style on column D).
$ws.Rows("330:331").Insert()

# --- New row 330 ---
$ws.Range("A330").Value() = 4
$ws.Range("B330").Value() = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C330").Value() = 'Los Lagos'
$ws.Range("D330").Value() = 44985
$ws.Range("E330").Value() = 10
$ws.Range("F330").Value() = 100112028
$ws.Range("G330").Value() = 'Sandia'
$ws.Range("H330").Value() = 'Sin especificar'
$ws.Range("I330").Value() = 'Primera'
$ws.Range("J330").Value() = 4000
$ws.Range("K330").Value() = 3000
$ws.Range("L330").Value() = 3200
$ws.Range("M330").Value() = 3075
$ws.Range("N330").Value() = '$/unidad'
$ws.Range("O330").Value() = "Región de O'Higgins"
$ws.Range("P330").Value() = 3075
$ws.Range("Q330").Value() = 1
$ws.Range("R330").Value() = 'Hortaliza'

# --- New row 331 ---
$ws.Range("A331").Value() = 4
$ws.Range("B331").Value() = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C331").Value() = 'Los Lagos'
$ws.Range("D331").Value() = 44985
$ws.Range("E331").Value() = 10
$ws.Range("F331").Value() = 100112028
$ws.Range("G331").Value() = 'Sandia'
$ws.Range("H331").Value() = 'Sin especificar'
$ws.Range("I331").Value() = 'Segunda'
$ws.Range("J331").Value() = 1500
$ws.Range("K331").Value() = 2500
$ws.Range("L331").Value() = 2500
$ws.Range("M331").Value() = 2500
$ws.Range("N331").Value() = '$/unidad'
$ws.Range("O331").Value() = "Región de O'Higgins"
$ws.Range("P331").Value() = 2500
$ws.Range("Q331").Value() = 1
$ws.Range("R331").Value() = 'Hortaliza'
